# Weekly update: insert a new "Jengibre" (ginger) price record as the new
# most-recent row for "Terminal La Palmera de La Serena", pushing the
# existing rows 32..116 down to 33..117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 32 (shifts 32..116 -> 33..117,
# carrying each row's values/styles down with it).
$ws.Rows.Item(32).EntireRow.Insert()

# Populate the newly-inserted row 32 with this week's data point.
$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 45012
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100114007
$ws.Range("G32").Value = "Jengibre"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 500
$ws.Range("K32").Value = 16000
$ws.Range("L32").Value = 17000
$ws.Range("M32").Value = 16500
$ws.Range("N32").Value = "$/caja 13 kilos"
$ws.Range("O32").Value = "Perú"
$ws.Range("P32").Value = 1269
$ws.Range("Q32").Value = 13
$ws.Range("R32").Value = "Hortaliza"
